# Insert a new data row before the existing row 72 (Mango / Macroferia
# Regional de Talca records). This shifts the former rows 72..200 down to
# 73..201 and grows the sheet's used range from A1:T200 to A1:T201.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with its record values.
$ws.Range("A72").Value = 5
$ws.Range("B72").Value = "Macroferia Regional de Talca"
$ws.Range("C72").Value = "Maule"
$ws.Range("D72").Value = 45210
$ws.Range("E72").Value = 7
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100108
$ws.Range("H72").Value = "Tropicales y subtropicales"
$ws.Range("I72").Value = 100108002
$ws.Range("J72").Value = "Mango"
$ws.Range("K72").Value = "Sin especificar"
$ws.Range("L72").Value = "Primera"
$ws.Range("M72").Value = 248
$ws.Range("N72").Value = 10000
$ws.Range("O72").Value = 10000
$ws.Range("P72").Value = 10000
$ws.Range("Q72").Value = "$/bandeja 4 kilos"
$ws.Range("R72").Value = "Brasil"
$ws.Range("S72").Value = 2500
$ws.Range("T72").Value = 4
